$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 360.83334
$ws.Range("I12").Value = 360.83334
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 360.83334
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -190.83334
$ws.Range("N12").ClearContents()
$ws.Range("H33").Value = 425.3
$ws.Range("I33").Value = 604.5
$ws.Range("J33").Value = 305.83334
$ws.Range("K33").Value = 604.5
$ws.Range("L33").Value = 305.83334
$ws.Range("M33").Value = -375.5
$ws.Range("N33").Value = -763.83334
$ws.Range("H40").Value = 3389.3333
$ws.Range("I40").Value = 3125.3333
$ws.Range("K40").Value = 3125.3333
$ws.Range("M40").Value = -2950.3333
$ws.Range("H80").Value = 1220.24
$ws.Range("I80").Value = 373.75
$ws.Range("J80").Value = 1618.5883
$ws.Range("K80").Value = 1121.25
$ws.Range("L80").Value = 4855.7649
$ws.Range("M80").Value = -123.25
$ws.Range("N80").Value = -6851.7649
$ws.Range("H83").Value = 1220.24
$ws.Range("I83").Value = 373.75
$ws.Range("J83").Value = 1618.5883
$ws.Range("K83").Value = 3363.75
$ws.Range("L83").Value = 14567.2947
$ws.Range("M83").Value = 1628.25
$ws.Range("N83").Value = -24551.2947
$ws.Range("H99").Value = 577.0833
$ws.Range("I99").Value = 372.6
$ws.Range("J99").Value = 1599.5
$ws.Range("K99").Value = 1117.8
$ws.Range("L99").Value = 4798.5
$ws.Range("M99").Value = 380.1999999999998
$ws.Range("N99").Value = -7794.5
$ws.Range("H127").Value = 1267.1538
$ws.Range("I127").Value = 1267.1538
$ws.Range("K127").Value = 3801.4614
$ws.Range("M127").Value = 1158.5386

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1879933.1
$ws.Range("I2").Value = 2046938.4
$ws.Range("J2").Value = 1125
$ws.Range("K2").Value = 2046938.4
$ws.Range("L2").Value = 1125
$ws.Range("M2").Value = -2046825.4
$ws.Range("N2").Value = -1351
$ws.Range("H5").Value = 324.5
$ws.Range("I5").Value = 400
$ws.Range("K5").Value = 400
$ws.Range("M5").Value = -288
$ws.Range("H32").Value = 57691.363
$ws.Range("I32").Value = 63642.95
$ws.Range("K32").Value = 63642.95
$ws.Range("M32").Value = -63355.95
$ws.Range("H45").Value = 150000
$ws.Range("J45").Value = 150000
$ws.Range("L45").Value = 150000
$ws.Range("N45").Value = -150754
$ws.Range("H116").Value = 1879933.1
$ws.Range("I116").Value = 2046938.4
$ws.Range("J116").Value = 1125
$ws.Range("K116").Value = 2046938.4
$ws.Range("L116").Value = 1125
$ws.Range("M116").Value = -2044644.4
$ws.Range("N116").Value = -5713
$ws.Range("H132").Value = 3857.8125
$ws.Range("I132").Value = 2887.4614
$ws.Range("K132").Value = 8662.3842
$ws.Range("M132").Value = -6132.3842

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1879933.1
$ws.Range("I3").Value = 2046938.4
$ws.Range("J3").Value = 1125
$ws.Range("K3").Value = 2046938.4
$ws.Range("L3").Value = 1125
$ws.Range("M3").Value = -2046824.4
$ws.Range("N3").Value = -1353
$ws.Range("H4").Value = 324.5
$ws.Range("I4").Value = 400
$ws.Range("K4").Value = 400
$ws.Range("M4").Value = -285
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -887
$ws.Range("H20").Value = 2473.074
$ws.Range("I20").Value = 2186.0557
$ws.Range("K20").Value = 2186.0557
$ws.Range("M20").Value = -1939.0557
$ws.Range("H22").Value = 5000.5
$ws.Range("I22").Value = 5000.5
$ws.Range("K22").Value = 5000.5
$ws.Range("M22").Value = -4827.5
$ws.Range("H86").Value = 1890.6923
$ws.Range("I86").Value = 1964.9166
$ws.Range("K86").Value = 1964.9166
$ws.Range("M86").Value = -841.9166
$ws.Range("H89").Value = 1890.6923
$ws.Range("I89").Value = 1964.9166
$ws.Range("K89").Value = 9824.583000000001
$ws.Range("M89").Value = -4208.583000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1939.75
$ws.Range("I16").Value = 1970.3334
$ws.Range("J16").Value = 1848
$ws.Range("K16").Value = 1970.3334
$ws.Range("L16").Value = 1848
$ws.Range("M16").Value = -1683.3334
$ws.Range("N16").Value = -2422
$ws.Range("H22").Value = 775
$ws.Range("I22").Value = 775
$ws.Range("K22").Value = 775
$ws.Range("M22").Value = -425
$ws.Range("H31").Value = 55560984
$ws.Range("J31").Value = 8181.4
$ws.Range("L31").Value = 8181.4
$ws.Range("N31").Value = -8771.4
$ws.Range("H34").Value = 55560984
$ws.Range("J34").Value = 8181.4
$ws.Range("L34").Value = 8181.4
$ws.Range("N34").Value = -8585.4
$ws.Range("H113").Value = 1939.75
$ws.Range("I113").Value = 1970.3334
$ws.Range("J113").Value = 1848
$ws.Range("K113").Value = 1970.3334
$ws.Range("L113").Value = 1848
$ws.Range("M113").Value = 199.6666
$ws.Range("N113").Value = -6188
$ws.Range("H122").Value = 113620.11
$ws.Range("I122").Value = 201516.2
$ws.Range("K122").Value = 604548.6000000001
$ws.Range("M122").Value = -602098.6000000001
$ws.Range("H134").Value = 4203.8
$ws.Range("I134").Value = 3079
$ws.Range("J134").Value = 6292.7144
$ws.Range("K134").Value = 9237
$ws.Range("L134").Value = 18878.1432
$ws.Range("M134").Value = -6702
$ws.Range("N134").Value = -23948.1432

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 586.75
$ws.Range("I22").Value = 415
$ws.Range("K22").Value = 1245
$ws.Range("M22").Value = -1076
$ws.Range("H27").Value = 586.75
$ws.Range("I27").Value = 415
$ws.Range("K27").Value = 1245
$ws.Range("M27").Value = -1143
$ws.Range("H104").Value = 3300
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("H106").Value = 8019.3335
$ws.Range("J106").Value = 8019.3335
$ws.Range("L106").Value = 24058.0005
$ws.Range("N106").Value = -25950.0005
$ws.Range("H107").Value = 1289.8182
$ws.Range("I107").Value = 205
$ws.Range("J107").Value = 1530.8889
$ws.Range("K107").Value = 615
$ws.Range("L107").Value = 4592.6667
$ws.Range("M107").Value = 1305
$ws.Range("N107").Value = -8432.6667
$ws.Range("H132").Value = 54557.473
$ws.Range("I132").Value = 72828
$ws.Range("K132").Value = 655452
$ws.Range("M132").Value = -652922

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5805.6665
$ws.Range("I70").Value = 4787.778
$ws.Range("K70").Value = 4787.778
$ws.Range("M70").Value = -4517.778
$ws.Range("H73").Value = 5805.6665
$ws.Range("I73").Value = 4787.778
$ws.Range("K73").Value = 4787.778
$ws.Range("M73").Value = -3851.778
$ws.Range("H80").Value = 6598.375
$ws.Range("J80").Value = 7131.5
$ws.Range("L80").Value = 7131.5
$ws.Range("N80").Value = -9127.5
$ws.Range("H83").Value = 6598.375
$ws.Range("J83").Value = 7131.5
$ws.Range("L83").Value = 35657.5
$ws.Range("N83").Value = -45641.5
$ws.Range("H126").Value = 3211.7
$ws.Range("I126").Value = 2167.2666
$ws.Range("J126").Value = 4256.1333
$ws.Range("K126").Value = 6501.7998
$ws.Range("L126").Value = 12768.3999
$ws.Range("M126").Value = -4031.7998
$ws.Range("N126").Value = -17708.3999
$ws.Range("H132").Value = 4472.5386
$ws.Range("I132").Value = 3498.0557
$ws.Range("K132").Value = 10494.1671
$ws.Range("M132").Value = -7964.167099999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 86940.836
$ws.Range("I7").Value = 127422.75
$ws.Range("J7").Value = 5977
$ws.Range("K7").Value = 127422.75
$ws.Range("L7").Value = 5977
$ws.Range("M7").Value = -127310.75
$ws.Range("N7").Value = -6201
$ws.Range("H46").Value = 6903.3794
$ws.Range("I46").Value = 1785.2858
$ws.Range("K46").Value = 1785.2858
$ws.Range("M46").Value = -1597.2858
$ws.Range("H55").Value = 144.85
$ws.Range("I55").Value = 107
$ws.Range("J55").Value = 182.7
$ws.Range("K55").Value = 107
$ws.Range("L55").Value = 182.7
$ws.Range("M55").Value = 66
$ws.Range("N55").Value = -528.7
$ws.Range("H122").Value = 2801
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 2502.5
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -12407.5
$ws.Range("H126").Value = 86940.836
$ws.Range("I126").Value = 127422.75
$ws.Range("J126").Value = 5977
$ws.Range("K126").Value = 382268.25
$ws.Range("L126").Value = 17931
$ws.Range("M126").Value = -379798.25
$ws.Range("N126").Value = -22871

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28273
$ws.Range("J41").Value = 28238.5
$ws.Range("L41").Value = 28238.5
$ws.Range("N41").Value = -29018.5
$ws.Range("H55").Value = 15196.429
$ws.Range("I55").Value = 3000
$ws.Range("K55").Value = 3000
$ws.Range("M55").Value = -2723
$ws.Range("H96").Value = 1975.15
$ws.Range("J96").Value = 2001
$ws.Range("L96").Value = 2001
$ws.Range("N96").Value = -4747
$ws.Range("H107").Value = 5471.2856
$ws.Range("I107").Value = 767.3333
$ws.Range("K107").Value = 2301.9999
$ws.Range("M107").Value = -381.9998999999998
$ws.Range("H126").Value = 5290.6113
$ws.Range("I126").Value = 4926.9375
$ws.Range("K126").Value = 14780.8125
$ws.Range("M126").Value = -12310.8125
$ws.Range("H131").Value = 132580.62
$ws.Range("J131").Value = 139999.28
$ws.Range("L131").Value = 139999.28
$ws.Range("N131").Value = -150079.28
